$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author inserted a new package entry ("bumpversion==0.5.3") above the
# existing list in column A. Column A is independent of the C/E columns
# (which already hold fill-down formulas for many rows), so the practical
# effect is just that the column-A values from row 7 down shift by one row,
# with the new entry landing in A7 and "Werkzeug==0.15.5" falling through
# to the newly used A25.
$values = @(
    "bumpversion==0.5.3",
    "Click==7.0",
    "et-xmlfile==1.0.1",
    "Flask==1.1.1",
    "Flask-Uploads==0.2.1",
    "itsdangerous==1.1.0",
    "jdcal==1.4.1",
    "Jinja2==2.10.1",
    "lml==0.0.9",
    "Markdown==3.1.1",
    "MarkupSafe==1.1.1",
    "openpyxl==2.5.14",
    "pip==18.1",
    "pkg-resources==0.0.0",
    "pyexcel-io==0.5.20",
    "pyexcel-xlsx==0.5.7",
    "setuptools==40.8.0",
    "SQLAlchemy==1.3.5",
    "Werkzeug==0.15.5"
)

$startRow = 7
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

# Restore the selection recorded in the saved workbook.
$ws.Range("M12").Select() | Out-Null
